# Apply updated crypto price/volume data per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.664.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "'1.635.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'212.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  +1.71%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("D10").Value = "'19.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.72%  "
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("D12").Value = "'1.864.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").Value = "'1.644.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").Value = "'4.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("E15").Value = "  +2.87%  "
$ws.Range("D16").Value = "'26.669.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "'63.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0₃0740"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").Value = "'209.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.43%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").Value = "'9.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("D23").Value = "'6.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").Value = "'1.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").Value = "'146.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("E28").Value = "  +2.98%  "
$ws.Range("D29").Value = "'15.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").Value = "'0.0520"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.14%  "
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("D36").Value = "'1.169.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("D38").Value = "'0.808"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.30%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "'0.504"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").Value = "'0.795"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").Value = "'1.776.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.10%  "
$ws.Range("D45").Value = "'92.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.0₆0104"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'54.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.409"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.30%  "
